# Update the last row (2025Q2, row 20) of the recurrence metrics sheet
# with the refreshed figures from the "bibi e add" dataset refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = 294
$ws.Range("D20").Value = 235
$ws.Range("E20").Value = 59
$ws.Range("F20").Value = 77.30263157894737
